# Update countries & provincias Spain
# Applies the 22-May-2020 11:05 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 11:05"

# Row 21 - Belgica: updated case counts (country unchanged)
$ws.Range("A21").Value = "Belgica"
$ws.Range("B21").Value = 56511
$ws.Range("C21").Value = 276
$ws.Range("D21").Value = 15123
$ws.Range("E21").Value = 32176
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 26
$ws.Range("H21").Value = 9212

# Row 30 - Banglades overtakes Portugal in the ranking
$ws.Range("A30").Value = "Banglades"
$ws.Range("B30").Value = 30205
$ws.Range("C30").Value = 1694
$ws.Range("D30").Value = 6190
$ws.Range("E30").Value = 23583
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 24
$ws.Range("H30").Value = 432

# Row 31 - Portugal moves down one place (same totals it had before)
$ws.Range("A31").Value = "Portugal"
$ws.Range("B31").Value = 29912
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 6452
$ws.Range("E31").Value = 22183
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 1277

# Row 34 - Indonesia overtakes Polonia in the ranking
$ws.Range("A34").Value = "Indonesia"
$ws.Range("B34").Value = 20796
$ws.Range("C34").Value = 634
$ws.Range("D34").Value = 5057
$ws.Range("E34").Value = 14413
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 48
$ws.Range("H34").Value = 1326

# Row 35 - Polonia moves down one place (same totals it had before)
$ws.Range("A35").Value = "Polonia"
$ws.Range("B35").Value = 20379
$ws.Range("C35").Value = 236
$ws.Range("D35").Value = 8731
$ws.Range("E35").Value = 10675
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 973

# Row 195 - Namibia overtakes San Vicente y las Granadinas in the ranking
$ws.Range("A195").Value = "Namibia"
$ws.Range("B195").Value = 19
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 14
$ws.Range("E195").Value = 5
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0

# Row 196 - San Vicente y las Granadinas moves down one place (same totals it had before)
$ws.Range("A196").Value = "San Vicente y las Granadinas"
$ws.Range("B196").Value = 18
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 14
$ws.Range("E196").Value = 4
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0
